$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "tyle" -> "tile" for the floor type in row 7
$ws.Range("A7").Value = "tile"

# Update the active selection to A8
$ws.Range("A8").Select()
